$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update model name labels in column A (rows reordered per new training run)
$ws.Range("A3").Value = "model_38_7_22"
$ws.Range("A4").Value = "model_38_7_21"
$ws.Range("A5").Value = "model_38_7_20"
$ws.Range("A6").Value = "model_38_7_19"
$ws.Range("A7").Value = "model_38_7_18"
$ws.Range("A8").Value = "model_38_7_17"
$ws.Range("A9").Value = "model_38_7_16"
$ws.Range("A10").Value = "model_38_7_15"
$ws.Range("A11").Value = "model_38_7_14"
$ws.Range("A12").Value = "model_38_7_13"
$ws.Range("A13").Value = "model_38_7_23"
$ws.Range("A15").Value = "model_38_7_10"
$ws.Range("A16").Value = "model_38_7_9"
$ws.Range("A17").Value = "model_38_7_8"
$ws.Range("A18").Value = "model_38_7_7"
$ws.Range("A19").Value = "model_38_7_6"
$ws.Range("A20").Value = "model_38_7_5"
$ws.Range("A21").Value = "model_38_7_4"
$ws.Range("A22").Value = "model_38_7_3"
$ws.Range("A23").Value = "model_38_7_2"
$ws.Range("A24").Value = "model_38_7_1"
$ws.Range("A25").Value = "model_38_7_11"

# Update metric columns B:Q - every row now shares the same aggregated metrics
$ws.Range("B2:B26").Value = 0.9999106709542923
$ws.Range("C2:C26").Value = 0.9989035467838921
$ws.Range("D2:D26").Value = 0.9998147383833075
$ws.Range("E2:E26").Value = 0.9997934851912618
$ws.Range("F2:F26").Value = 0.9998341903495032
$ws.Range("G2:G26").Value = 0.00008338474307615185
$ws.Range("H2:H26").Value = 0.001023490948502973
$ws.Range("I2:I26").Value = 0.0002206355631893498
$ws.Range("J2:J26").Value = 0.0001571851393505272
$ws.Range("K2:K26").Value = 0.0001889103512242937
$ws.Range("L2:L26").Value = 0.0005467859721176896
$ws.Range("M2:M26").Value = 0.009131524685185483
$ws.Range("N2:N26").Value = 1.000064966578696
$ws.Range("O2:O26").Value = 0.009520272595895827
$ws.Range("P2:P26").Value = 132.7840904041427
$ws.Range("Q2:Q26").Value = 202.2600124216301
